# The second sheet ("Resp Responde pelo cumprimento ...") has a column B
# of evidence ratings where several labels were entered with the space
# between "Muito/Pouco/Não" and "Evidenciado" missing (a typo duplicated
# against the correctly-spaced labels already used elsewhere in the same
# column, e.g. B2/B7/B12 etc. use "Extremamente Evidenciado"). Re-typing
# the affected cells with the correctly spaced text both fixes the label
# and lets the workbook collapse the now-unused duplicate shared-string
# entries away on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

$ws.Range("B3").Value = "Muito Evidenciado"
$ws.Range("B5").Value = "Pouco Evidenciado"
$ws.Range("B6").Value = "Não Evidenciado"
$ws.Range("B8").Value = "Muito Evidenciado"
$ws.Range("B10").Value = "Pouco Evidenciado"
$ws.Range("B11").Value = "Não Evidenciado "
$ws.Range("B13").Value = "Muito Evidenciado"
$ws.Range("B15").Value = "Pouco Evidenciado"
$ws.Range("B16").Value = "Não Evidenciado"
$ws.Range("B18").Value = "Muito Evidenciado"
$ws.Range("B20").Value = "Pouco Evidenciado"
$ws.Range("B21").Value = "Não Evidenciado"
$ws.Range("B23").Value = "Muito Evidenciado"
$ws.Range("B25").Value = "Pouco Evidenciado"
$ws.Range("B26").Value = "Não Evidenciado"

# Leave the cursor where the author's saved view shows it (top-left A1,
# selection on B26) instead of the previous D7.
$ws.Range("A1").Select()
$ws.Range("B26").Select()
